# Auto-generated edit script applying the Valefor_Profits.xlsx diff
# Updates FFXIV market-price snapshot columns (H-N) across the 8 server sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 168.6
$ws.Range("I12").Value = 148
$ws.Range("J12").Value = 216.66667
$ws.Range("K12").Value = 148
$ws.Range("L12").Value = 216.66667
$ws.Range("M12").Value = 22
$ws.Range("N12").Value = -556.6666700000001
$ws.Range("H33").Value = 30310782
$ws.Range("I33").Value = 33335174
$ws.Range("J33").Value = 66856.664
$ws.Range("K33").Value = 33335174
$ws.Range("L33").Value = 66856.664
$ws.Range("M33").Value = -33334945
$ws.Range("H40").Value = 1648.8572
$ws.Range("I40").Value = 1836.421
$ws.Range("J40").Value = 1493.9131
$ws.Range("K40").Value = 1836.421
$ws.Range("L40").Value = 1493.9131
$ws.Range("M40").Value = -1661.421
$ws.Range("N40").Value = -1843.9131
$ws.Range("H86").Value = 3572.5
$ws.Range("I86").Value = 5114.4287
$ws.Range("J86").Value = 2591.2727
$ws.Range("K86").Value = 5114.4287
$ws.Range("L86").Value = 2591.2727
$ws.Range("M86").Value = -3991.4287
$ws.Range("N86").Value = -4837.2727
$ws.Range("H89").Value = 3572.5
$ws.Range("I89").Value = 5114.4287
$ws.Range("J89").Value = 2591.2727
$ws.Range("K89").Value = 25572.1435
$ws.Range("L89").Value = 12956.3635
$ws.Range("M89").Value = -19956.1435
$ws.Range("N89").Value = -24188.3635
$ws.Range("H98").Value = 25205.977
$ws.Range("I98").Value = 28127.842
$ws.Range("J98").Value = 2999.8
$ws.Range("K98").Value = 28127.842
$ws.Range("L98").Value = 2999.8
$ws.Range("M98").Value = -26629.842
$ws.Range("H122").Value = 25205.977
$ws.Range("I122").Value = 28127.842
$ws.Range("J122").Value = 2999.8
$ws.Range("K122").Value = 84383.526
$ws.Range("L122").Value = 8999.400000000001
$ws.Range("M122").Value = -81933.526
$ws.Range("H132").Value = 2233726
$ws.Range("I132").Value = 2791096.5
$ws.Range("J132").Value = 4244.4165
$ws.Range("K132").Value = 8373289.5
$ws.Range("L132").Value = 12733.2495
$ws.Range("M132").Value = -8370759.5
$ws.Range("N132").Value = -17793.2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2968.2222
$ws.Range("I35").Value = 2669
$ws.Range("J35").Value = 3566.6667
$ws.Range("K35").Value = 2669
$ws.Range("L35").Value = 3566.6667
$ws.Range("M35").Value = -2263
$ws.Range("N35").Value = -4378.6667
$ws.Range("H41").Value = 2802.5557
$ws.Range("I41").Value = 2379.1765
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 2379.1765
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -1965.1765
$ws.Range("N41").Value = -10828
$ws.Range("H61").Value = 3252.5
$ws.Range("I61").Value = 2909.4119
$ws.Range("J61").Value = 4085.7144
$ws.Range("K61").Value = 2909.4119
$ws.Range("L61").Value = 4085.7144
$ws.Range("M61").Value = -2697.4119
$ws.Range("N61").Value = -4509.7144
$ws.Range("H110").Value = 1261.8064
$ws.Range("I110").Value = 698.6842
$ws.Range("J110").Value = 2153.4167
$ws.Range("K110").Value = 698.6842
$ws.Range("L110").Value = 2153.4167
$ws.Range("M110").Value = 1346.3158
$ws.Range("N110").Value = -6243.4167
$ws.Range("H136").Value = 3252.5
$ws.Range("I136").Value = 2909.4119
$ws.Range("J136").Value = 4085.7144
$ws.Range("K136").Value = 8728.235700000001
$ws.Range("L136").Value = 12257.1432
$ws.Range("M136").Value = -6178.235700000001
$ws.Range("N136").Value = -17357.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 6930
$ws.Range("I37").Value = 600
$ws.Range("J37").Value = 8512.5
$ws.Range("K37").Value = 600
$ws.Range("L37").Value = 8512.5
$ws.Range("M37").Value = -463
$ws.Range("N37").Value = -8786.5
$ws.Range("H45").Value = 40355
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 40355
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 40355
$ws.Range("N45").Value = -41971
$ws.Range("H54").Value = 2157.75
$ws.Range("I54").Value = 793.6667
$ws.Range("J54").Value = 6250
$ws.Range("K54").Value = 793.6667
$ws.Range("L54").Value = 6250
$ws.Range("M54").Value = -309.6667
$ws.Range("N54").Value = -7218

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 810.5714
$ws.Range("I22").Value = 533.1667
$ws.Range("J22").Value = 1018.625
$ws.Range("K22").Value = 533.1667
$ws.Range("L22").Value = 1018.625
$ws.Range("M22").Value = -183.1667
$ws.Range("N22").Value = -1718.625
$ws.Range("H31").Value = 15876138
$ws.Range("I31").Value = 25642200
$ws.Range("J31").Value = 6287.8335
$ws.Range("K31").Value = 25642200
$ws.Range("L31").Value = 6287.8335
$ws.Range("M31").Value = -25641905
$ws.Range("N31").Value = -6877.8335
$ws.Range("H34").Value = 15876138
$ws.Range("I34").Value = 25642200
$ws.Range("J34").Value = 6287.8335
$ws.Range("K34").Value = 25642200
$ws.Range("L34").Value = 6287.8335
$ws.Range("M34").Value = -25641998
$ws.Range("N34").Value = -6691.8335
$ws.Range("H64").Value = 24613.334
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 24613.334
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 24613.334
$ws.Range("N64").Value = -25109.334
$ws.Range("H67").Value = 24613.334
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 24613.334
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 24613.334
$ws.Range("N67").Value = -26329.334
$ws.Range("H141").Value = 31584.615
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 31584.615
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 31584.615
$ws.Range("N141").Value = -41944.61500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 610.2222
$ws.Range("I5").Value = 417.11765
$ws.Range("J5").Value = 938.5
$ws.Range("K5").Value = 1251.35295
$ws.Range("L5").Value = 2815.5
$ws.Range("M5").Value = -1139.35295
$ws.Range("N5").Value = -3039.5
$ws.Range("H50").Value = 410.91666
$ws.Range("I50").Value = 268.33334
$ws.Range("J50").Value = 458.44446
$ws.Range("K50").Value = 805.0000200000001
$ws.Range("L50").Value = 1375.33338
$ws.Range("M50").Value = -324.0000200000001
$ws.Range("N50").Value = -2337.33338
$ws.Range("H53").Value = 410.91666
$ws.Range("I53").Value = 268.33334
$ws.Range("J53").Value = 458.44446
$ws.Range("K53").Value = 805.0000200000001
$ws.Range("L53").Value = 1375.33338
$ws.Range("M53").Value = -324.0000200000001
$ws.Range("N53").Value = -2337.33338
$ws.Range("H113").Value = 553.1667
$ws.Range("I113").Value = 546
$ws.Range("J113").Value = 555.55554
$ws.Range("K113").Value = 1638
$ws.Range("L113").Value = 1666.66662
$ws.Range("M113").Value = 532
$ws.Range("N113").Value = -6006.66662
$ws.Range("H135").Value = 610.2222
$ws.Range("I135").Value = 417.11765
$ws.Range("J135").Value = 938.5
$ws.Range("K135").Value = 3754.05885
$ws.Range("L135").Value = 8446.5
$ws.Range("M135").Value = -1219.05885
$ws.Range("N135").Value = -13516.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1999.3889
$ws.Range("I107").Value = 2349.0908
$ws.Range("J107").Value = 1449.8572
$ws.Range("K107").Value = 2349.0908
$ws.Range("L107").Value = 1449.8572
$ws.Range("M107").Value = -429.0907999999999
$ws.Range("N107").Value = -5289.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1179.4546
$ws.Range("I46").Value = 890
$ws.Range("J46").Value = 1288
$ws.Range("K46").Value = 890
$ws.Range("L46").Value = 1288
$ws.Range("M46").Value = -702
$ws.Range("N46").Value = -1664
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -14002
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 45000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -40008
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H132").Value = 2485.238
$ws.Range("I132").Value = 1636.9375
$ws.Range("J132").Value = 5199.8
$ws.Range("K132").Value = 4910.8125
$ws.Range("L132").Value = 15599.4
$ws.Range("M132").Value = -2380.8125
$ws.Range("N132").Value = -20659.4
$ws.Range("H136").Value = 3604.9443
$ws.Range("I136").Value = 3149.2144
$ws.Range("J136").Value = 5200
$ws.Range("K136").Value = 9447.643199999999
$ws.Range("L136").Value = 15600
$ws.Range("M136").Value = -6897.643199999999
$ws.Range("N136").Value = -20700
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1520.7018
$ws.Range("I132").Value = 1262.8372
$ws.Range("J132").Value = 2312.7144
$ws.Range("K132").Value = 3788.5116
$ws.Range("L132").Value = 6938.1432
$ws.Range("M132").Value = -1258.5116
$ws.Range("N132").Value = -11998.1432
$ws.Range("H136").Value = 2281.8845
$ws.Range("I136").Value = 1828.909
$ws.Range("J136").Value = 2614.0667
$ws.Range("K136").Value = 5486.727000000001
$ws.Range("L136").Value = 7842.2001
$ws.Range("M136").Value = -2936.727000000001
$ws.Range("N136").Value = -12942.2001

